$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Body sectPr: drop the explicit footnote-restart-per-section setting so
#    the section falls back to the (default) continuous footnote numbering.
# ---------------------------------------------------------------------------
$d.Footnotes.NumberingRule = 0   # wdRestartContinuous

# ---------------------------------------------------------------------------
# 2) Title / TitleChar styles: tighten the big title font (adds character
#    spacing + kerning threshold) used for the paper title.
# ---------------------------------------------------------------------------
$titleStyle = $d.Styles("Title")
$titleStyle.Font.Spacing = -0.5   # -10 twentieths-of-a-point
$titleStyle.Font.Kerning = 14     # 28 half-points

$titleCharStyle = $d.Styles("TitleChar")
$titleCharStyle.Font.Spacing = -0.5
$titleCharStyle.Font.Kerning = 14

# ---------------------------------------------------------------------------
# 3) Author / Date styles: stop inheriting from Title (so they no longer
#    pick up the large title font) and instead explicitly center them.
# ---------------------------------------------------------------------------
$authorStyle = $d.Styles("Author")
$authorStyle.BaseStyle = $null
$authorStyle.ParagraphFormat.Alignment = 1   # wdAlignParagraphCenter

$dateStyle = $d.Styles("Date")
$dateStyle.BaseStyle = $null
$dateStyle.ParagraphFormat.Alignment = 1   # wdAlignParagraphCenter
